# Update the "Metadata" sheet with the new Title, Date and Description values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Title (row 5, column B) - new value
$ws.Range("B5").Value = "NG-Imm Client HIV Status VS"

# Date (row 8, column B) - updated timestamp
$ws.Range("B8").Value = "2025-06-24T09:13:37+01:00"

# Description (row 13, column B) - unchanged text value, rewritten so it
# becomes its own shared-string entry now that the old one was repurposed.
$ws.Range("B13").Value = "Client HIV Status"
